# Update price list values per commit "Cambio en lista de precios"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 33-39 (VC Balanced / VC Complete dog food)
$ws.Range("B33").Value = 3610.0
$ws.Range("B34").Value = 3180.0
$ws.Range("B35").Value = 3180.0
$ws.Range("B36").Value = 4470.0
$ws.Range("B37").Value = 2590.0
$ws.Range("B38").Value = 2470.0
$ws.Range("B39").Value = 2350.0

# Rows 66-71 (VC Gato Balanced / VC Gato Complete cat food)
$ws.Range("B66").Value = 7390.0
$ws.Range("B67").Value = 6690.0
$ws.Range("B68").Value = 6580.0
$ws.Range("B69").Value = 7470.0
$ws.Range("B70").Value = 4080.0
$ws.Range("B71").Value = 3925.0
